# Add "car functions" support to the MapUsers workbook:
#   - a new "صالح الماضي" / IRQ200 user row (row 9)
#   - a new "carcode" column (F) on Table1, populated per-row
#   - refreshed selection / dimension bookkeeping

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# 1) Append a new row to the table for the new user "صالح الماضي" (IRQ200)
$newRow = $lo.ListRows.Add()
$ws.Range("A9").Value = "صالح الماضي"
$ws.Range("B9").Value = "IRQ200"
$ws.Range("C9").Value = "صالح الماضي"
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 105

# 2) Append a new "carcode" column to the table
$newCol = $lo.ListColumns.Add()
$ws.Range("F1").Value = "carcode"

# 3) Populate the carcode values for the existing rows
$ws.Range("F7").Value = "waritex11123"
$ws.Range("F4").Value = "waritex107"
$ws.Range("F3").Value = "waritex3"
$ws.Range("F6").Value = "waritex6"
$ws.Range("F5").Value = "waritex7"
$ws.Range("F2").Value = "waritex8"
# rows 8 and 9 have no carcode yet (left blank)

# 4) Match formatting used by the rest of the table (center/center alignment)
$ws.Range("E1").Copy()
$ws.Range("A9:F9").PasteSpecial(-4122)
$ws.Range("F1:F8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 5) Resize the new column to fit its contents and refresh the selection
$ws.Range("F1:F9").EntireColumn.AutoFit() | Out-Null
$null = $ws.Range("F10").Select()
